# Updates the "Price" (D) and "Volume(1h)" (E) columns on the crypto
# tracker sheet to the latest scraped snapshot values.
#
# Column D holds plain text in the source data (prices such as
# "35.528.60" that use '.' as both a thousands- and decimal-separator,
# so they are never genuinely numeric) but some of the individual
# values (e.g. "0.710", "247.09") happen to *look* like valid numbers.
# A bare `Range.Value = "..."` assignment lets Excel's normal type
# inference silently convert those into floating point numbers, which
# would corrupt the text and introduce float rounding noise. To keep
# such values as literal text (matching the original file) we briefly
# flip the cell to the "Text" number format before assigning the
# value, then restore its original (General) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $row, $col, $value) {
    $looksNumeric = $value -match '^[+-]?[0-9]*\.?[0-9]+$'

    $cell = $ws.Cells.Item($row, $col)
    if ($looksNumeric) {
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = $origStyle
    } else {
        $cell.Value = $value
    }
}

$updates = @(
    @{ Row = 2; D = "35.528.60"; E = "  +0.04%  " },
    @{ Row = 3; D = "1.913.40"; E = "  +0.32%  " },
    @{ Row = 4; D = $null; E = "  -0.15%  " },
    @{ Row = 5; D = "0.710"; E = "  +9.15%  " },
    @{ Row = 6; D = "247.09"; E = "  +0.21%  " },
    @{ Row = 7; D = $null; E = "  -0.06%  " },
    @{ Row = 8; D = "40.79"; E = "  -3.12%  " },
    @{ Row = 9; D = $null; E = "  +4.36%  " },
    @{ Row = 10; D = "52.69"; E = "  +7.54%  " },
    @{ Row = 11; D = $null; E = "  +2.56%  " },
    @{ Row = 12; D = $null; E = "  -1.01%  " },
    @{ Row = 13; D = "2.189.91"; E = "  +0.31%  " },
    @{ Row = 14; D = "12.70"; E = "  +2.25%  " },
    @{ Row = 15; D = "0.718"; E = "  +2.46%  " },
    @{ Row = 16; D = "1.917.58"; E = "  +0.54%  " },
    @{ Row = 17; D = $null; E = "  +1.41%  " },
    @{ Row = 18; D = "35.515.16"; E = "  -0.03%  " },
    @{ Row = 19; D = "73.41"; E = "  +1.64%  " },
    @{ Row = 20; D = "0.0₃0829"; E = "  -0.85%  " },
    @{ Row = 21; D = "13.14"; E = "  +3.84%  " },
    @{ Row = 22; D = "242.50"; E = "  -0.56%  " },
    @{ Row = 23; D = "5.06"; E = "  +4.23%  " },
    @{ Row = 24; D = $null; E = "  -0.07%  " },
    @{ Row = 25; D = "2.33"; E = "  +0.97%  " },
    @{ Row = 26; D = "2.31"; E = "  +4.53%  " },
    @{ Row = 27; D = "168.74"; E = "  -1.52%  " },
    @{ Row = 28; D = "8.67"; E = "  +1.52%  " },
    @{ Row = 29; D = "18.80"; E = "  +4.08%  " },
    @{ Row = 30; D = $null; E = "  +4.45%  " },
    @{ Row = 31; D = "4.119.75"; E = $null },
    @{ Row = 32; D = "4.25"; E = "  +1.78%  " },
    @{ Row = 33; D = "0.0578"; E = "  +1.25%  " },
    @{ Row = 34; D = $null; E = "  +0.19%  " },
    @{ Row = 35; D = $null; E = "  +7.59%  " },
    @{ Row = 36; D = $null; E = "  -0.13%  " },
    @{ Row = 37; D = "0.919"; E = "  -5.24%  " },
    @{ Row = 38; D = "1.48"; E = "  +11.15%  " },
    @{ Row = 39; D = "2.05"; E = "  +0.63%  " },
    @{ Row = 40; D = "17.34"; E = "  +10.47%  " },
    @{ Row = 41; D = "98.25"; E = "  +5.81%  " },
    @{ Row = 42; D = "1.14"; E = "  +3.14%  " },
    @{ Row = 43; D = $null; E = "  +2.53%  " },
    @{ Row = 44; D = "0.0648"; E = "  +1.48%  " },
    @{ Row = 45; D = "1.355.63"; E = "  +0.49%  " },
    @{ Row = 46; D = "2.46"; E = "  +2.56%  " },
    @{ Row = 47; D = $null; E = "  +0.33%  " },
    @{ Row = 48; D = "2.79"; E = "  +1.11%  " },
    @{ Row = 49; D = $null; E = "  -5.78%  " },
    @{ Row = 50; D = "12.26"; E = "  -3.71%  " },
    @{ Row = 51; D = "6.57"; E = "  -0.53%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-TextCell $ws $u.Row 4 $u.D
    }
    if ($null -ne $u.E) {
        Set-TextCell $ws $u.Row 5 $u.E
    }
}
